$wb = $excel.ActiveWorkbook

# --- instructions sheet: add new rows 13-17 with notes ---
$ws1 = $wb.Worksheets.Item("instructions")

$ws1.Range("A13").Value = "12-bit Addresses/12-bit offsets"
$ws1.Range("A14").Value = "r0-r7 are general purpose"
$ws1.Range("A15").Value = "r8 = 0x00"
$ws1.Range("A16").Value = "r9 = 0x01"
$ws1.Range("A17").Value = "r15 = 0xFFr"
$ws1.Range("A13:A17").HorizontalAlignment = -4131

# --- zoom each sheet to 220% ---
foreach ($name in @("instructions","example","MemMap")) {
    $s = $wb.Worksheets.Item($name)
    $s.Activate() | Out-Null
    $excel.ActiveWindow.Zoom = 220
}

# --- restore per-sheet selections ---
$wb.Worksheets.Item("example").Range("D5").Select() | Out-Null
$wb.Worksheets.Item("MemMap").Range("H5").Select() | Out-Null

# --- make "instructions" the active sheet with A14 selected ---
$ws1.Activate() | Out-Null
$ws1.Range("A14").Select() | Out-Null
